# Add 2022-Q4 data
# 1) Insert a new "2022-Q4" sheet right after "总计" (cloned from "2022-Q2" so it
#    inherits identical sheet formatting / page setup), populate it with the
#    four new fund rows.
# 2) Insert a new summary row into "总计" for the 2022-Q4 quarter, shifting the
#    existing quarter rows down by one.

$wb = $excel.ActiveWorkbook

function Set-TextNoStyle($ws, $addr, $val, $formatSourceAddr) {
    # Writing a numeric-looking string via .Value normally gets coerced to a
    # number by Excel. Force text storage by switching the cell to the Text
    # number format before assignment, then borrow the (unstyled) format of
    # a neighbouring cell so we don't leave a stray "Text" style behind.
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $ws.Range($formatSourceAddr).Copy()
    $cell.PasteSpecial(-4122)  # xlPasteFormats
}

# ---------------------------------------------------------------------------
# Step 1: new "2022-Q4" sheet
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q2Sheet.Copy($null, $totalSheet)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q4"

# The template only has 2 data rows; stamp its row-2 formatting (incl. the
# styled A column) down onto rows 4 and 5 so every data row matches.
$newSheet.Range("A2:H2").Copy()
$newSheet.Range("A4").PasteSpecial(-4122)
$newSheet.Range("A5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$neutral = "H2"  # an already-unstyled cell to borrow "no style" formatting from

$newSheet.Range("A2").Value = 0
Set-TextNoStyle $newSheet "B2" "004871" $neutral
Set-TextNoStyle $newSheet "C2" "中银金融地产混合A" $neutral
Set-TextNoStyle $newSheet "D2" "1.58" $neutral
Set-TextNoStyle $newSheet "E2" "89.99" $neutral
Set-TextNoStyle $newSheet "F2" "2.94" $neutral
Set-TextNoStyle $newSheet "G2" "0.0465" $neutral
$newSheet.Range("H2").Value = 9

$newSheet.Range("A3").Value = 1
Set-TextNoStyle $newSheet "B3" "010312" $neutral
Set-TextNoStyle $newSheet "C3" "中银金融地产混合C" $neutral
Set-TextNoStyle $newSheet "D3" "0.51" $neutral
Set-TextNoStyle $newSheet "E3" "89.99" $neutral
Set-TextNoStyle $newSheet "F3" "2.94" $neutral
Set-TextNoStyle $newSheet "G3" "0.0150" $neutral
$newSheet.Range("H3").Value = 9

$newSheet.Range("A4").Value = 2
Set-TextNoStyle $newSheet "B4" "011494" $neutral
Set-TextNoStyle $newSheet "C4" "华泰紫金丰和偏债混合发起A" $neutral
Set-TextNoStyle $newSheet "D4" "0.12" $neutral
Set-TextNoStyle $newSheet "E4" "36.70" $neutral
Set-TextNoStyle $newSheet "F4" "1.03" $neutral
Set-TextNoStyle $newSheet "G4" "0.0012" $neutral
$newSheet.Range("H4").Value = 9

$newSheet.Range("A5").Value = 3
Set-TextNoStyle $newSheet "B5" "011495" $neutral
Set-TextNoStyle $newSheet "C5" "华泰紫金丰和偏债混合发起C" $neutral
Set-TextNoStyle $newSheet "D5" "0.03" $neutral
Set-TextNoStyle $newSheet "E5" "36.70" $neutral
Set-TextNoStyle $newSheet "F5" "1.03" $neutral
Set-TextNoStyle $newSheet "G5" "0.0003" $neutral
$newSheet.Range("H5").Value = 9

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Step 2: insert the 2022-Q4 row into "总计", shifting old rows down
# ---------------------------------------------------------------------------
for ($r = 6; $r -ge 2; $r--) {
    $dest = $r + 1
    $totalSheet.Cells.Item($dest, 1).Value = $dest - 2
    $totalSheet.Cells.Item($dest, 2).Value = $totalSheet.Cells.Item($r, 2).Value()
    $totalSheet.Cells.Item($dest, 3).Value = $totalSheet.Cells.Item($r, 3).Value()
    $totalSheet.Cells.Item($dest, 4).Value = $totalSheet.Cells.Item($r, 4).Value()
}

# Row 7 is brand-new (didn't exist before the shift), so it needs column A's
# style (s="2") copied across manually; row 2-6 already carry it forward.
$totalSheet.Range("A6").Copy()
$totalSheet.Range("A7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q4"
$totalSheet.Cells.Item(2, 3).Value = 4
$totalSheet.Cells.Item(2, 4).Value = 0.06

# ---------------------------------------------------------------------------
# Step 3: restore the originally-active tab ("2021-Q1", now the last sheet)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
